$wb = $excel.ActiveWorkbook

# --- 1. Update the Date value on the Metadata sheet ---
$wsMeta = $wb.Worksheets.Item("Metadata")
$wsMeta.Range("B8").Value = "2022-05-03T12:45:05-04:00"

# --- 2. Replace the LOINC-include sheet contents with an
#        "Include from Observation Valu 2" style CodeSystem include,
#        pointing at the new SPLASCHFrequencyCS code system ---
$ws3 = $wb.Worksheets.Item("Include from LOINC")

# Drop the old per-code answer rows (rows 5-8): Often/Always plus the
# blank separator + System URI/loinc.org rows - they get rebuilt below.
$ws3.Range("A5:B8").EntireRow.Delete()

# Row 1 header becomes a single "Codes" cell (no second column)
$ws3.Range("A1").Value = "Codes"
$ws3.Range("B1").Clear()

# Row 2 becomes a single "All codes" cell (no second column)
$ws3.Range("A2").Value = "All codes"
$ws3.Range("B2").Clear()

# Row 3 stays a blank separator row (already blank/empty in the template)
$ws3.Range("A3").Value = ""
$ws3.Range("B3").Value = ""

# Row 4 becomes the System URI pointing at the new CodeSystem
$ws3.Range("A4").Value = "System URI"
$ws3.Range("B4").Value = "http://hl7.org/fhir/us/pacio-splasch/CodeSystem/SPLASCHFrequencyCS"

# --- 3. Rename the sheet tab to match its new content ---
$ws3.Name = "Include from Observation Valu 2"
